# 20140513 4 new solved
# Move four previously-unsolved interview questions (rows that contained
# question #30/#29/#31/#8) up to the top of the "unsolved" block (rows 42-45),
# mark them solved (column A = 1) and add the remark text that explains how
# each one was solved (column D). The remaining still-unsolved rows shift
# down to fill rows 46-49, keeping their original (blank) A/D values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Newly solved rows (moved to the top of the block, rows 42-45)
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "面试题30"
$ws.Range("C42").Value = "最小的K个数"
$ws.Range("D42").Value = "微软面试，nth element，也可最小堆"
$ws.Range("E42").Value = "已收录"

$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "面试题29"
$ws.Range("C43").Value = "数组中出现次数超过一半的数字"
$ws.Range("D43").Value = "动态维护指针和数组"
$ws.Range("E43").Value = "已收录"

$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "面试题31"
$ws.Range("C44").Value = "连续子数组的最大和"
$ws.Range("D44").Value = "动态规划，和小于零重新开始"
$ws.Range("E44").Value = "已收录"

$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "面试题8"
$ws.Range("C45").Value = "旋转数组的最小数字"
$ws.Range("D45").Value = "二分查找，相等只能顺序"
$ws.Range("E45").Value = "已收录"

# Remaining still-unsolved rows, shifted down into rows 46-49 (unchanged
# content, still blank A/D)
$ws.Range("A46").Value = ""
$ws.Range("B46").Value = "面试题1"
$ws.Range("C46").Value = "赋值运算符函数"
$ws.Range("D46").Value = ""
$ws.Range("E46").Value = "不适合在线模式"

$ws.Range("A47").Value = ""
$ws.Range("B47").Value = "面试题2"
$ws.Range("C47").Value = "实现Singleton模式"
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = "不适合在线模式"

$ws.Range("A48").Value = ""
$ws.Range("B48").Value = "面试题18"
$ws.Range("C48").Value = "树的子结构"
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = "已收录"

$ws.Range("A49").Value = ""
$ws.Range("B49").Value = "面试题20"
$ws.Range("C49").Value = "顺时针打印矩阵"
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = "已收录"

$wb.Save()
